$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: paragraph 1 - append a red "(This is a change ... )" note
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$r = $p1.Range
$r.End = $r.End - 1          # exclude the paragraph mark
$r.Collapse(0)                # wdCollapseEnd
$r.InsertAfter("  ")
$r.Collapse(0)

$enDash = [char]0x2013
$chunk1 = "(This is a change " + $enDash + " Ve"
$chunk2 = "rsion for branch alternate"
$chunk3 = ")"

$start1 = $r.Start
$r.InsertAfter($chunk1)
$d.Range($start1, $start1 + $chunk1.Length).Font.Color = 192
$r.Collapse(0)

$start2 = $r.Start
$r.InsertAfter($chunk2)
$d.Range($start2, $start2 + $chunk2.Length).Font.Color = 192
$r.Collapse(0)

$start3 = $r.Start
$r.InsertAfter($chunk3)
$d.Range($start3, $start3 + $chunk3.Length).Font.Color = 192
$r.Collapse(0)

# ---------------------------------------------------------------------
# Change 2: "Crispian's Day speech from Shakespear's Henry V [Source..]"
# paragraph - tidy up run splitting around "from " / "Shakespear's" and
# merge the trailing "Henry V [Source - Wikipedia]" runs into one.
# ---------------------------------------------------------------------
$p4 = $d.Paragraphs.Item(4)
$pStart = $p4.Range.Start
$pText = $p4.Range.Text

$idxFrom = $pText.IndexOf("from")
$fromEnd = $pStart + $idxFrom + 4
$idxShakespear = $pText.IndexOf("Shakespear")
$shakespearStart = $pStart + $idxShakespear
$gapRange = $d.Range($fromEnd, $shakespearStart)
$gapRange.Delete()
$fromRange = $d.Range($pStart + $idxFrom, $fromEnd)
$fromRange.InsertAfter(" ")

$p4b = $d.Paragraphs.Item(4)
$pStart2 = $p4b.Range.Start
$pText2 = $p4b.Range.Text

$idxHenry = $pText2.IndexOf("Henry")
$henryStart = $pStart2 + $idxHenry - 1   # include the leading space
$idxCloseBracket = $pText2.IndexOf("]")
$closeBracketEnd = $pStart2 + $idxCloseBracket + 1

$idxDash = $pText2.IndexOf([char]0x2013)
$dashStart = $pStart2 + $idxDash
$dashEnd = $dashStart + 1

$afterDash = $d.Range($dashEnd, $closeBracketEnd)
$afterDash.Delete()
$beforeDash = $d.Range($henryStart, $dashStart)
$beforeDash.Delete()

$mergedRange = $d.Range($henryStart, $henryStart + 1)
$mergedRange.Text = " Henry V [Source " + $enDash + " Wikipedia]"

# ---------------------------------------------------------------------
# Change 3: append two trailing paragraphs before the end of the body
#           (a "larger"-styled empty paragraph, then a bare paragraph)
# ---------------------------------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0)   # wdCollapseEnd
$wordml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p><w:pPr><w:pStyle w:val="larger"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="150" w:afterAutospacing="0"/></w:pPr></w:p>
<w:p/>
</w:body></w:document>
</pkg:xmlData></pkg:part></pkg:package>
'@
$endRange.InsertXML($wordml)

# ---------------------------------------------------------------------
# Change 4: styles.xml clean-up - drop unused "apple-converted-space"
#           and "Hyperlink" character styles (delete higher index first
#           so the lower one's lookup-by-name stays valid).
# ---------------------------------------------------------------------
$d.Styles.Item("Hyperlink").Delete()
$d.Styles.Item("apple-converted-space").Delete()
